# Insert a new data row at row 135 (weekly update for Fruta/Hortaliza),
# shifting the existing rows 135:152 down to 136:153.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("135").Insert()

$ws.Cells.Item(135, 1).Value  = 1
$ws.Cells.Item(135, 2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(135, 3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(135, 4).Value  = 44522
$ws.Cells.Item(135, 5).Value  = 15
$ws.Cells.Item(135, 6).Value  = 'Fruta'
$ws.Cells.Item(135, 7).Value  = 100108
$ws.Cells.Item(135, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(135, 9).Value  = 100108006
$ws.Cells.Item(135, 10).Value = 'Plátano'
$ws.Cells.Item(135, 11).Value = 'Sin especificar'
$ws.Cells.Item(135, 12).Value = 'Maduro'
$ws.Cells.Item(135, 13).Value = 120
$ws.Cells.Item(135, 14).Value = 12000
$ws.Cells.Item(135, 15).Value = 13000
$ws.Cells.Item(135, 16).Value = 12500
$ws.Cells.Item(135, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(135, 18).Value = 'Ecuador'
$ws.Cells.Item(135, 19).Value = 625
$ws.Cells.Item(135, 20).Value = 20
